# Auto-generated script to apply cryptos.xlsx cell value updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.408.62'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.282.76'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '498.21'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.36'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('E9').Value = '  +3.39%  '
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.333'
$ws.Range('E11').Value = '  +3.68%  '
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').Value = '2.689.08'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.45'
$ws.Range('E14').Value = '  +5.04%  '
$ws.Range('D15').Value = '54.341.65'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '2.284.85'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.25'
$ws.Range('E18').Value = '  +6.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.14'
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '305.83'
$ws.Range('E20').Value = '  +2.70%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.44'
$ws.Range('E21').Value = '  +2.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '62.12'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '2.382.57'
$ws.Range('E25').Value = '  +2.91%  '
$ws.Range('E26').Value = '  +2.79%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.32'
$ws.Range('E27').Value = '  +3.08%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '173.67'
$ws.Range('E28').Value = '  +6.65%  '
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('D30').Value = '0.0₃0690'
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.95'
$ws.Range('E31').Value = '  +2.37%  '
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.81'
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.926'
$ws.Range('E36').Value = '  +10.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.21'
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.77'
$ws.Range('E38').Value = '  +4.84%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.374'
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.42'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '126.26'
$ws.Range('E42').Value = '  -1.27%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.80'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0492'
$ws.Range('E44').Value = '  +3.68%  '
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.549'
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '240.44'
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('B48').Value = 'Polygon'
$ws.Range('C48').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.373'
$ws.Range('E48').Value = '  +0.84%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0206'
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.77'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '16.33'
$ws.Range('E51').Value = '  +0.87%  '
